# Re-order data rows 2..27 according to the permutation observed between
# the "before" and "after" versions of the workbook. Each destination row
# ends up with the *entire* contents (every column A:AY) that used to live
# in a different source row; row 1 (headers) and rows 7, 11 and 18 keep
# their own original content (they map to themselves).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map: destination row number -> source row number (1-based worksheet rows)
$rowMap = @{
    2  = 4
    3  = 8
    4  = 3
    5  = 22
    6  = 24
    7  = 7
    8  = 23
    9  = 26
    10 = 15
    11 = 11
    12 = 20
    13 = 2
    14 = 21
    15 = 13
    16 = 19
    17 = 16
    18 = 18
    19 = 5
    20 = 6
    21 = 14
    22 = 27
    23 = 17
    24 = 25
    25 = 9
    26 = 12
    27 = 10
}

$firstCol = 1   # A
$lastCol  = 51  # AY
$firstRow = 2
$lastRow  = 27
$numCols  = $lastCol - $firstCol + 1

# Snapshot every source row's values (columns A:AY) before writing anything,
# since several rows are both a source and a destination in the permutation.
$snapshot = @{}
for ($r = $firstRow; $r -le $lastRow; $r++) {
    $rowValues = New-Object 'object[]' $numCols
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $rowValues[$c - $firstCol] = $ws.Cells.Item($r, $c).Value2
    }
    $snapshot[$r] = $rowValues
}

# Write the permuted rows back out, but only touch a cell when its new
# value actually differs from what is already there. This avoids Excel's
# automatic "smart" re-typing (e.g. turning a date-looking text value into
# a real date) for cells whose content is not actually changing.
for ($destRow = $firstRow; $destRow -le $lastRow; $destRow++) {
    $srcRow = $rowMap[$destRow]
    if ($srcRow -eq $destRow) {
        continue
    }
    $rowValues = $snapshot[$srcRow]
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $newVal = $rowValues[$c - $firstCol]
        $cell = $ws.Cells.Item($destRow, $c)
        $curVal = $cell.Value2
        if ($curVal -ne $newVal) {
            $cell.Value = $newVal
        }
    }
}
